$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Generik" column header in D1, matching the header style
# used by the existing header cells (A1:C1).
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("D1").Value = "Generik"
